$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.8680022954940796
$ws.Range("B1").Value = 1.443885803222656
$ws.Range("C1").Value = 3.891918182373047
$ws.Range("D1").Value = 2.615260601043701
$ws.Range("E1").Value = 1.583036780357361
